$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to Text so that numeric-looking strings
# (e.g. "19.80", "449.00") are not silently coerced into numbers, which
# would both lose formatting (trailing zeros) and change the stored cell type.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "68.031.78"
$ws.Range("E2").Value = "  +1.56%  "

# Row 3
$ws.Range("D3").Value = "3.901.76"
$ws.Range("E3").Value = "  -0.73%  "

# Row 4
$ws.Range("E4").Value = "  -0.12%  "

# Row 5
$ws.Range("D5").Value = "483.76"
$ws.Range("E5").Value = "  +2.45%  "

# Row 6
$ws.Range("D6").Value = "145.42"
$ws.Range("E6").Value = "  +0.43%  "

# Row 7
$ws.Range("D7").Value = "0.625"
$ws.Range("E7").Value = "  -1.22%  "

# Row 8
$ws.Range("E8").Value = "  -0.13%  "

# Row 9
$ws.Range("D9").Value = "0.725"
$ws.Range("E9").Value = "  -3.56%  "

# Row 10
$ws.Range("D10").Value = "0.167"
$ws.Range("E10").Value = "  +2.96%  "

# Row 11
$ws.Range("D11").Value = "0.0000362"
$ws.Range("E11").Value = "  +11.76%  "

# Row 12
$ws.Range("D12").Value = "42.76"
$ws.Range("E12").Value = "  -1.60%  "

# Row 13
$ws.Range("D13").Value = "10.63"
$ws.Range("E13").Value = "  +0.45%  "

# Row 14
$ws.Range("D14").Value = "4.515.83"
$ws.Range("E14").Value = "  -0.78%  "

# Row 15
$ws.Range("D15").Value = "14.73"
$ws.Range("E15").Value = "  -2.31%  "

# Row 16
$ws.Range("D16").Value = "3.899.20"
$ws.Range("E16").Value = "  -1.72%  "

# Row 17
$ws.Range("E17").Value = "  -0.41%  "

# Row 18
$ws.Range("D18").Value = "19.80"
$ws.Range("E18").Value = "  -2.37%  "

# Row 19
$ws.Range("E19").Value = "  -4.14%  "

# Row 20
$ws.Range("D20").Value = "68.079.53"
$ws.Range("E20").Value = "  +0.83%  "

# Row 21
$ws.Range("D21").Value = "449.00"
$ws.Range("E21").Value = "  +3.75%  "

# Row 22
$ws.Range("E22").Value = "  -2.36%  "

# Row 23
$ws.Range("D23").Value = "3.35"
$ws.Range("E23").Value = "  -0.61%  "

# Row 24
$ws.Range("D24").Value = "89.58"
$ws.Range("E24").Value = "  +0.84%  "

# Row 25
$ws.Range("D25").Value = "11.62"
$ws.Range("E25").Value = "  +12.83%  "

# Row 26
$ws.Range("B26").Value = "RenderToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D26").Value = "10.63"
$ws.Range("E26").Value = "  +8.02%  "

# Row 27
$ws.Range("D27").Value = "38.81"
$ws.Range("E27").Value = "  -1.29%  "

# Row 28
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").Value = "3.57"
$ws.Range("E28").Value = "  -0.76%  "

# Row 29
$ws.Range("E29").Value = "  +3.02%  "

# Row 30
$ws.Range("D30").Value = "13.38"
$ws.Range("E30").Value = "  -2.99%  "

# Row 31
$ws.Range("D31").Value = "689.57"
$ws.Range("E31").Value = "  -6.30%  "

# Row 32
$ws.Range("E32").Value = "  -2.40%  "

# Row 33
$ws.Range("D33").Value = "2.84"
$ws.Range("E33").Value = "  +1.78%  "

# Row 34
$ws.Range("D34").Value = "0.0₃0975"
$ws.Range("E34").Value = "  +34.07%  "

# Row 35
$ws.Range("E35").Value = "  -6.32%  "

# Row 36
$ws.Range("D36").Value = "59.20"
$ws.Range("E36").Value = "  +0.99%  "

# Row 37
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").Value = "5.66"
$ws.Range("E37").Value = "  +4.62%  "

# Row 38
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "0.149"
$ws.Range("E38").Value = "  -7.92%  "

# Row 39
$ws.Range("D39").Value = "0.998"
$ws.Range("E39").Value = "  -0.18%  "

# Row 40
$ws.Range("D40").Value = "0.0473"
$ws.Range("E40").Value = "  -2.47%  "

# Row 41
$ws.Range("D41").Value = "2.77"
$ws.Range("E41").Value = "  +8.59%  "

# Row 42
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "3.05"
$ws.Range("E42").Value = "  +8.25%  "

# Row 43
$ws.Range("B43").Value = "ThetaToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D43").Value = "3.01"
$ws.Range("E43").Value = "  -2.90%  "

# Row 44
$ws.Range("D44").Value = "0.347"
$ws.Range("E44").Value = "  +1.50%  "

# Row 45
$ws.Range("E45").Value = "  -0.48%  "

# Row 46
$ws.Range("E46").Value = "  -0.18%  "

# Row 47
$ws.Range("E47").Value = "  -2.32%  "

# Row 48
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").Value = "147.16"
$ws.Range("E48").Value = "  +2.36%  "

# Row 49
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").Value = "2.12"
$ws.Range("E49").Value = "  -4.81%  "

# Row 50
$ws.Range("E50").Value = "  -3.01%  "

# Row 51
$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D51").Value = "0.000264"
$ws.Range("E51").Value = "  +71.39%  "

# Remove the explicit Text style we applied above so the cells revert to
# the workbooks original unstyled state (value type stays Text).
$ws.Range("D2:E51").ClearFormats()